$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the header row (row 1) with new labels
$ws.Range("B1").Value = "Question"
$ws.Range("A1").Value = "Question Type"
$ws.Range("C1").Value = "Correct"
$ws.Range("D1").Value = "Options"
$ws.Range("E1").Value = "Marks"
$ws.Range("F1").Value = "Question For"

# Update the "Question For" value in row 2
$ws.Range("F2").Value = "natural sciences, education, philosophy"

# Update selection
$ws.Range("F4").Select()
